$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33 (ALC)
$ws.Range("H33").Value = 266778.2
$ws.Range("I33").Value = 34755.277
$ws.Range("J33").Value = 1948944.2
$ws.Range("K33").Value = 34755.277
$ws.Range("L33").Value = 1948944.2
$ws.Range("M33").Value = -34526.277
$ws.Range("N33").Value = -1949402.2

# Row 41 (ALC)
$ws.Range("H41").Value = 742.02856
$ws.Range("I41").Value = 1293.0769
$ws.Range("J41").Value = 416.4091
$ws.Range("K41").Value = 1293.0769
$ws.Range("L41").Value = 416.4091
$ws.Range("M41").Value = -853.0769
$ws.Range("N41").Value = -1296.4091

# Row 51 (ALC)
$ws.Range("H51").Value = 8050
$ws.Range("I51").Value = 12840.1
$ws.Range("J51").Value = 3259.9
$ws.Range("K51").Value = 12840.1
$ws.Range("L51").Value = 3259.9
$ws.Range("M51").Value = -12356.1
$ws.Range("N51").Value = -4227.9

# Row 62 (ALC)
$ws.Range("H62").Value = 2090.4546
$ws.Range("J62").Value = 1698.3334
$ws.Range("L62").Value = 1698.3334
$ws.Range("N62").Value = -2946.3334

# Row 65 (ALC)
$ws.Range("H65").Value = 2090.4546
$ws.Range("J65").Value = 1698.3334
$ws.Range("L65").Value = 8491.666999999999
$ws.Range("N65").Value = -14731.667

# Row 93 (ALC)
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

# Row 100 (ALC)
$ws.Range("H100").Value = 957
$ws.Range("I100").Value = 781
$ws.Range("J100").Value = 1103.6666
$ws.Range("K100").Value = 781
$ws.Range("L100").Value = 1103.6666
$ws.Range("M100").Value = -240
$ws.Range("N100").Value = -2185.6666

# Row 128 (ALC)
$ws.Range("H128").Value = 37596.668
$ws.Range("I128").Value = 43000
$ws.Range("J128").Value = 34895
$ws.Range("K128").Value = 43000
$ws.Range("L128").Value = 34895
$ws.Range("M128").Value = -38020
$ws.Range("N128").Value = -44855

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 3039.32
$ws.Range("I32").Value = 2683.2827
$ws.Range("J32").Value = 7133.75
$ws.Range("K32").Value = 2683.2827
$ws.Range("L32").Value = 7133.75
$ws.Range("M32").Value = -2396.2827
$ws.Range("N32").Value = -7707.75

# Row 110 (ARM)
$ws.Range("H110").Value = 125250750
$ws.Range("I110").Value = 143143460
$ws.Range("K110").Value = 143143460
$ws.Range("M110").Value = -143141415

# Row 132 (ARM)
$ws.Range("H132").Value = 1918.9114
$ws.Range("I132").Value = 1945.9155
$ws.Range("J132").Value = 1679.25
$ws.Range("K132").Value = 5837.7465
$ws.Range("L132").Value = 5037.75
$ws.Range("M132").Value = -3307.7465
$ws.Range("N132").Value = -10097.75

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (BSM)
$ws.Range("H20").Value = 23218.174
$ws.Range("I20").Value = 33908.516
$ws.Range("J20").Value = 1124.8
$ws.Range("K20").Value = 33908.516
$ws.Range("L20").Value = 1124.8
$ws.Range("M20").Value = -33661.516
$ws.Range("N20").Value = -1618.8

# Row 134 (BSM)
$ws.Range("H134").Value = 1956.6721
$ws.Range("I134").Value = 1776.5283
$ws.Range("J134").Value = 3150.125
$ws.Range("K134").Value = 5329.5849
$ws.Range("L134").Value = 9450.375
$ws.Range("M134").Value = -2794.5849
$ws.Range("N134").Value = -14520.375

$ws = $wb.Worksheets.Item("CRP")
# Row 43 (CRP)
$ws.Range("H43").Value = 21444.334
$ws.Range("J43").Value = 21444.334
$ws.Range("L43").Value = 21444.334
$ws.Range("N43").Value = -21812.334

# Row 99 (CRP)
$ws.Range("H99").Value = 20298.334
$ws.Range("I99").Value = 6060
$ws.Range("J99").Value = 34536.668
$ws.Range("K99").Value = 6060
$ws.Range("L99").Value = 34536.668
$ws.Range("M99").Value = -4562
$ws.Range("N99").Value = -37532.668

# Row 101 (CRP)
$ws.Range("H101").Value = 21444.334
$ws.Range("J101").Value = 21444.334
$ws.Range("L101").Value = 21444.334
$ws.Range("N101").Value = -27934.334

# Row 106 (CRP)
$ws.Range("H106").Value = 25223.666
$ws.Range("J106").Value = 25223.666
$ws.Range("L106").Value = 25223.666
$ws.Range("N106").Value = -27747.666

# Row 126 (CRP)
$ws.Range("H126").Value = 20298.334
$ws.Range("I126").Value = 6060
$ws.Range("J126").Value = 34536.668
$ws.Range("K126").Value = 18180
$ws.Range("L126").Value = 103610.004
$ws.Range("M126").Value = -15710
$ws.Range("N126").Value = -108550.004

$ws = $wb.Worksheets.Item("CUL")
# Row 15 (CUL)
$ws.Range("H15").Value = 81.052635
$ws.Range("I15").Value = 42.727272
$ws.Range("J15").Value = 133.75
$ws.Range("K15").Value = 128.181816
$ws.Range("L15").Value = 401.25
$ws.Range("M15").Value = 11.818184
$ws.Range("N15").Value = -681.25

# Row 109 (CUL)
$ws.Range("H109").Value = 3662.7727
$ws.Range("I109").Value = 2390.875
$ws.Range("J109").Value = 4389.5713
$ws.Range("K109").Value = 7172.625
$ws.Range("L109").Value = 13168.7139
$ws.Range("M109").Value = -6132.625
$ws.Range("N109").Value = -15248.7139

$ws = $wb.Worksheets.Item("GSM")
# Row 5 (GSM)
$ws.Range("H5").Value = 5000000
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

# Row 70 (GSM)
$ws.Range("H70").Value = 73682.92999999999
$ws.Range("I70").Value = 121855.06
$ws.Range("J70").Value = 5439.0835
$ws.Range("K70").Value = 121855.06
$ws.Range("L70").Value = 5439.0835
$ws.Range("M70").Value = -121585.06
$ws.Range("N70").Value = -5979.0835

# Row 73 (GSM)
$ws.Range("H73").Value = 73682.92999999999
$ws.Range("I73").Value = 121855.06
$ws.Range("J73").Value = 5439.0835
$ws.Range("K73").Value = 121855.06
$ws.Range("L73").Value = 5439.0835
$ws.Range("M73").Value = -120919.06
$ws.Range("N73").Value = -7311.0835

# Row 126 (GSM)
$ws.Range("H126").Value = 8406202
$ws.Range("I126").Value = 3358.3333
$ws.Range("J126").Value = 14708336
$ws.Range("K126").Value = 10074.9999
$ws.Range("L126").Value = 44125008
$ws.Range("M126").Value = -7604.999899999999
$ws.Range("N126").Value = -44129948

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (LTW)
$ws.Range("H22").Value = 1384.3462
$ws.Range("J22").Value = 1424.75
$ws.Range("L22").Value = 1424.75
$ws.Range("N22").Value = -2014.75

# Row 27 (LTW)
$ws.Range("H27").Value = 1384.3462
$ws.Range("J27").Value = 1424.75
$ws.Range("L27").Value = 1424.75
$ws.Range("N27").Value = -1638.75

# Row 40 (LTW)
$ws.Range("H40").Value = 34170.484
$ws.Range("I40").Value = 101498.5
$ws.Range("J40").Value = 2109.524
$ws.Range("K40").Value = 101498.5
$ws.Range("L40").Value = 2109.524
$ws.Range("M40").Value = -101362.5
$ws.Range("N40").Value = -2381.524

# Row 61 (LTW)
$ws.Range("H61").Value = 1409.0769
$ws.Range("I61").Value = 1231.8334
$ws.Range("J61").Value = 1807.875
$ws.Range("K61").Value = 1231.8334
$ws.Range("L61").Value = 1807.875
$ws.Range("M61").Value = -1029.8334
$ws.Range("N61").Value = -2211.875

# Row 111 (LTW)
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

# Row 113 (LTW)
$ws.Range("H113").Value = 1409.0769
$ws.Range("I113").Value = 1231.8334
$ws.Range("J113").Value = 1807.875
$ws.Range("K113").Value = 1231.8334
$ws.Range("L113").Value = 1807.875
$ws.Range("M113").Value = 938.1666
$ws.Range("N113").Value = -6147.875

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (WVR)
$ws.Range("H81").Value = 251717.38
$ws.Range("I81").Value = 334266.66
$ws.Range("J81").Value = 202187.8
$ws.Range("K81").Value = 668533.3199999999
$ws.Range("L81").Value = 404375.6
$ws.Range("M81").Value = -667472.3199999999
$ws.Range("N81").Value = -406497.6

# Row 84 (WVR)
$ws.Range("H84").Value = 251717.38
$ws.Range("I84").Value = 334266.66
$ws.Range("J84").Value = 202187.8
$ws.Range("K84").Value = 3342666.6
$ws.Range("L84").Value = 2021878
$ws.Range("M84").Value = -3337362.6
$ws.Range("N84").Value = -2032486

# Row 113 (WVR)
$ws.Range("H113").Value = 645.9655
$ws.Range("I113").Value = 506.3684
$ws.Range("J113").Value = 911.2
$ws.Range("K113").Value = 1519.1052
$ws.Range("L113").Value = 2733.6
$ws.Range("M113").Value = 650.8948
$ws.Range("N113").Value = -7073.6
